# ---------------------------------------------------------------------------
# Add PPTX presentations and publication linking
#
# 1. "Projects" sheet: rename the "slides" column to "presentations" and
#    rewrite its per-project entries in the new "Title::file.pptx::image url"
#    shape (dropping the old slide-title/url/caption triples). The separate
#    "cover" image column (D) is cleared out - projects no longer carry a
#    standalone cover image/hyperlink.
# 2. Add a new "Publications" sheet listing papers/articles with links back
#    to the Projects and Researchers sheets via slugs.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Projects sheet updates
# ---------------------------------------------------------------------------
$projects = $wb.Worksheets.Item("Projects")

# Drop the "cover" column's data + hyperlinks (column header D1 "cover"
# itself is left in place, only the per-project values go away).
$projects.Hyperlinks.Delete()
$projects.Range("D2:D4").ClearContents()
$projects.Range("D2:D4").Style = "Normal"

# "slides" -> "presentations", with a new compact pptx-file based format.
$projects.Range("G1").Value = "presentations"
$projects.Range("G2").Value = "Predict Overview::predict-overview.pptx::https://placehold.co/960x540/png?text=Predict+Deck+1|Predict Update::predict-update.pptx::https://placehold.co/960x540/png?text=Predict+Deck+2"
$projects.Range("G3").Value = "Demand Forecast Brief::nursing-demand-brief.pptx::https://placehold.co/960x540/png?text=Demand+Deck+1"
$projects.Range("G4").Value = "Attrition Insights::attrition-insights.pptx::https://placehold.co/960x540/png?text=Attrition+Deck+1"

# ---------------------------------------------------------------------------
# 2. New "Publications" sheet (added after the last existing sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$pubs = $wb.Worksheets.Add($null, $lastSheet)
$pubs.Name = "Publications"

$headers = @("slug", "title", "url", "authors", "year", "project", "researchers")
for ($col = 1; $col -le $headers.Count; $col++) {
    $pubs.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# slug, title, url, authors, year, project slug (blank if none), researcher slugs
$rows = @(
    @("operational-ai-in-practice", "Operational AI in Practice", "https://example.com/operational-ai", "Theofilos, Diego Gonzalez Garcia-Torres", "2024", "predict-operational-ai", "theofilos|diego-gonzalez-garcia-torres"),
    @("forecasting-nurse-demand", "Forecasting Nurse Demand", "https://example.com/nurse-demand", "Maya Singh, Jordan Lee", "2025", "nursing-demand-forecast", "maya-singh|jordan-lee"),
    @("attrition-phenotypes-at-scale", "Attrition Phenotypes at Scale", "https://example.com/attrition-phenotypes", "Maya Singh, Theofilos", "2025", "attrition-phenotyping", "maya-singh|theofilos"),
    @("operational-ai-decision-support", "Operational AI Decision Support", "https://example.com/decision-support", "Theofilos, Jordan Lee", "2024", "predict-operational-ai", "theofilos|jordan-lee"),
    @("wellbeing-signals-in-staffing", "Well-being Signals in Staffing", "https://example.com/wellbeing-signals", "Maya Singh, Diego Gonzalez Garcia-Torres", "2023", "", "maya-singh|diego-gonzalez-garcia-torres"),
    @("systemwide-workforce-analytics", "System-wide Workforce Analytics", "https://example.com/workforce-analytics", "Theofilos, Maya Singh", "2023", "", "theofilos|maya-singh")
)

# Years read as text (shared strings), matching the rest of the sheet's data.
$pubs.Range("E2:E7").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $pubs.Cells.Item($r, 1).Value = $row[0]
    $pubs.Cells.Item($r, 2).Value = $row[1]
    $pubs.Cells.Item($r, 3).Value = $row[2]
    $pubs.Cells.Item($r, 4).Value = $row[3]
    $pubs.Cells.Item($r, 5).Value = $row[4]
    if ($row[5] -ne "") {
        $pubs.Cells.Item($r, 6).Value = $row[5]
    }
    $pubs.Cells.Item($r, 7).Value = $row[6]
    $r++
}

# Hyperlink the url column (C) for every publication row, same pattern as
# the existing sheets use for their external links.
for ($r = 2; $r -le 7; $r++) {
    $cell = $pubs.Cells.Item($r, 3)
    $pubs.Hyperlinks.Add($cell, $cell.Value) | Out-Null
}
